$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-06-15 Saturday" "2024-06-16 Sunday"

Replace-Text "407×4=1628" "681×9=6129"
Replace-Text "140×3=420" "171×3=513"
Replace-Text "746×9=6714" "924×7=6468"
Replace-Text "882×4=3528" "981×8=7848"
Replace-Text "769×2=1538" "115×6=690"

Replace-Text "243×7=1701" "301×9=2709"
Replace-Text "520×7=3640" "518×7=3626"
Replace-Text "365×8=2920" "333×5=1665"
Replace-Text "653×2=1306" "245×8=1960"
Replace-Text "987×8=7896" "464×4=1856"

Replace-Text "159×8=1272" "368×3=1104"
Replace-Text "761×3=2283" "582×3=1746"
Replace-Text "690×5=3450" "404×5=2020"
Replace-Text "647×2=1294" "313×4=1252"
Replace-Text "139×4=556" "577×4=2308"

Replace-Text "835×4=3340" "709×5=3545"
Replace-Text "545×6=3270" "374×7=2618"
Replace-Text "765×9=6885" "450×5=2250"
Replace-Text "969×8=7752" "585×2=1170"
Replace-Text "837×9=7533" "608×9=5472"

Replace-Text "911×4=3644" "863×6=5178"
Replace-Text "704×2=1408" "863×3=2589"
Replace-Text "385×6=2310" "393×6=2358"
Replace-Text "114×8=912" "858×3=2574"
Replace-Text "125×5=625" "992×3=2976"
